# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (linked from the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours      (linked from the Slide Master / Design)
#
# The authored edit swaps the two themes' colour schemes, so the slides
# (which follow theme2 via the Slide Master / Design) end up using the
# plain "Office Theme" palette instead of "Integral".
#
# PowerPoint's automation model only exposes the *active* design's theme
# colours for editing (Presentation.SlideMaster.Theme.ThemeColorScheme),
# so we reassign each of its 12 theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) to the values that the "Office Theme"
# palette used, matching the authored diff.

$p = $ppt.ActivePresentation
$design = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme palette (was previously only used by ppt/theme/theme1.xml)
# RGB() isn't available in this host, so colours are passed as the packed
# 0xBBGGRR long value PowerPoint stores on ColorFormat.RGB.
$officeThemeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $design.Colors($i).RGB = $officeThemeColors[$i - 1]
}
